$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values for rows 4-8, columns A,B,D,E,F,G,H
$cols = @("A","B","D","E","F","G","H")
$orig = @{}
foreach ($r in 4..8) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Mapping of source row -> destination row, based on the diff:
# before row 4 -> after row 6
# before row 5 -> after row 4
# before row 6 -> after row 5
# before row 7 -> after row 8
# before row 8 -> after row 7
$map = @{ 4 = 6; 5 = 4; 6 = 5; 7 = 8; 8 = 7 }

foreach ($srcRow in $map.Keys) {
    $dstRow = $map[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$dstRow").Value2 = $orig[$srcRow][$c]
    }
}
